# Bump the cached "datetimeFigureOut" footer date (4/15/2022 -> 4/19/2022)
# everywhere it is stored: the slide master and every slide layout's Date
# placeholder. Slides themselves inherit the footer from the
# master/layout, so none of the individual slides need touching.

$p = $ppt.ActivePresentation

$oldDate = "4/15/2022"
$newDate = "4/19/2022"

function Update-DatePlaceholder($container) {
    if ($container.Shapes -eq $null) { return }
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shape = $container.Shapes.Item($i)

        $phType = $null
        try { $phType = $shape.PlaceholderFormat.Type } catch {}

        if ($phType -eq 16) {
            if ($shape.HasTextFrame) {
                $tr = $shape.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}
